$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 814.8103599999999
$ws.Range("J17").Value = 814.8103599999999
$ws.Range("L17").Value = 2444.43108
$ws.Range("N17").Value = -2780.43108
$ws.Range("H18").Value = 41669424
$ws.Range("I18").Value = 47620936
$ws.Range("K18").Value = 47620936
$ws.Range("M18").Value = -47620652
$ws.Range("H136").Value = 137912
$ws.Range("J136").Value = 137912
$ws.Range("L136").Value = 137912
$ws.Range("N136").Value = -148112
$ws.Range("H139").Value = 130333
$ws.Range("J139").Value = 130333
$ws.Range("L139").Value = 130333
$ws.Range("N139").Value = -140613
$ws.Range("H140").Value = 59783.8
$ws.Range("J140").Value = 58569.89
$ws.Range("L140").Value = 58569.89
$ws.Range("N140").Value = -68929.89

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 490.83334
$ws.Range("I5").Value = 490.83334
$ws.Range("K5").Value = 490.83334
$ws.Range("M5").Value = -378.83334
$ws.Range("H25").Value = 649.5
$ws.Range("I25").Value = 649.5
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 649.5
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -247.5
$ws.Range("N25").ClearContents()
$ws.Range("H32").Value = 6435.927
$ws.Range("I32").Value = 3738.1
$ws.Range("J32").Value = 13793.637
$ws.Range("K32").Value = 3738.1
$ws.Range("L32").Value = 13793.637
$ws.Range("M32").Value = -3451.1
$ws.Range("N32").Value = -14367.637

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 490.83334
$ws.Range("I4").Value = 490.83334
$ws.Range("K4").Value = 490.83334
$ws.Range("M4").Value = -375.83334
$ws.Range("H86").Value = 1882
$ws.Range("I86").Value = 1936.5714
$ws.Range("K86").Value = 1936.5714
$ws.Range("M86").Value = -813.5714
$ws.Range("H89").Value = 1882
$ws.Range("I89").Value = 1936.5714
$ws.Range("K89").Value = 9682.857
$ws.Range("M89").Value = -4066.857

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1007.1429
$ws.Range("I22").Value = 283.33334
$ws.Range("K22").Value = 283.33334
$ws.Range("M22").Value = 66.66665999999998
$ws.Range("H26").Value = 4000
$ws.Range("I26").Value = 4000
$ws.Range("K26").Value = 4000
$ws.Range("M26").Value = -3713
$ws.Range("H41").Value = 104000
$ws.Range("I41").Value = 47500
$ws.Range("J41").Value = 141666.67
$ws.Range("K41").Value = 47500
$ws.Range("L41").Value = 141666.67
$ws.Range("M41").Value = -47072
$ws.Range("N41").Value = -142522.67
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H58").Value = 4525.125
$ws.Range("I58").Value = 4507.4287
$ws.Range("J58").Value = 4649
$ws.Range("K58").Value = 4507.4287
$ws.Range("L58").Value = 4649
$ws.Range("M58").Value = -4304.4287
$ws.Range("N58").Value = -5055
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H86").Value = 5750
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 5750
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H107").Value = 503.58334
$ws.Range("I107").Value = 429.3
$ws.Range("K107").Value = 429.3
$ws.Range("M107").Value = 1490.7
$ws.Range("H109").Value = 69984.5
$ws.Range("J109").Value = 69984.5
$ws.Range("L109").Value = 69984.5
$ws.Range("N109").Value = -72064.5
$ws.Range("H115").Value = 53000
$ws.Range("J115").Value = 53000
$ws.Range("L115").Value = 53000
$ws.Range("N115").Value = -55350
$ws.Range("H136").Value = 4525.125
$ws.Range("I136").Value = 4507.4287
$ws.Range("J136").Value = 4649
$ws.Range("K136").Value = 13522.2861
$ws.Range("L136").Value = 13947
$ws.Range("M136").Value = -10972.2861
$ws.Range("N136").Value = -19047
$ws.Range("H141").Value = 96432.08
$ws.Range("J141").Value = 96432.08
$ws.Range("L141").Value = 96432.08
$ws.Range("N141").Value = -106792.08

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 922
$ws.Range("I98").Value = 593
$ws.Range("J98").Value = 947.3077
$ws.Range("K98").Value = 1779
$ws.Range("L98").Value = 2841.9231
$ws.Range("M98").Value = -281
$ws.Range("N98").Value = -5837.9231
$ws.Range("H131").Value = 8623305
$ws.Range("J131").Value = 7249139
$ws.Range("L131").Value = 21747417
$ws.Range("N131").Value = -21757497

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 7574321
$ws.Range("I14").Value = 17667100
$ws.Range("J14").Value = 4737.25
$ws.Range("K14").Value = 17667100
$ws.Range("L14").Value = 4737.25
$ws.Range("M14").Value = -17666932
$ws.Range("N14").Value = -5073.25
$ws.Range("H80").Value = 80299.06
$ws.Range("I80").Value = 171966.67
$ws.Range("J80").Value = 25298.5
$ws.Range("K80").Value = 171966.67
$ws.Range("L80").Value = 25298.5
$ws.Range("M80").Value = -170968.67
$ws.Range("N80").Value = -27294.5
$ws.Range("H83").Value = 80299.06
$ws.Range("I83").Value = 171966.67
$ws.Range("J83").Value = 25298.5
$ws.Range("K83").Value = 859833.3500000001
$ws.Range("L83").Value = 126492.5
$ws.Range("M83").Value = -854841.3500000001
$ws.Range("N83").Value = -136476.5
$ws.Range("H102").Value = 21747390
$ws.Range("I102").Value = 35723428
$ws.Range("J102").Value = 6888.3335
$ws.Range("K102").Value = 35723428
$ws.Range("L102").Value = 6888.3335
$ws.Range("M102").Value = -35721806
$ws.Range("N102").Value = -10132.3335
$ws.Range("H140").Value = 78739.5
$ws.Range("J140").Value = 78739.5
$ws.Range("L140").Value = 78739.5
$ws.Range("N140").Value = -89099.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4867.846
$ws.Range("I7").Value = 1569.5
$ws.Range("J7").Value = 6333.778
$ws.Range("K7").Value = 1569.5
$ws.Range("L7").Value = 6333.778
$ws.Range("M7").Value = -1457.5
$ws.Range("N7").Value = -6557.778
$ws.Range("H46").Value = 7044.6
$ws.Range("J46").Value = 7044.6
$ws.Range("L46").Value = 7044.6
$ws.Range("N46").Value = -7420.6
$ws.Range("H109").Value = 99000
$ws.Range("J109").Value = 99000
$ws.Range("L109").Value = 99000
$ws.Range("N109").Value = -101774
$ws.Range("H122").Value = 49696776
$ws.Range("I122").Value = 71432056
$ws.Range("J122").Value = 15886342
$ws.Range("K122").Value = 214296168
$ws.Range("L122").Value = 47659026
$ws.Range("M122").Value = -214293718
$ws.Range("N122").Value = -47663926
$ws.Range("H126").Value = 4867.846
$ws.Range("I126").Value = 1569.5
$ws.Range("J126").Value = 6333.778
$ws.Range("K126").Value = 4708.5
$ws.Range("L126").Value = 19001.334
$ws.Range("M126").Value = -2238.5
$ws.Range("N126").Value = -23941.334
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H131").Value = 100000
$ws.Range("J131").Value = 100000
$ws.Range("L131").Value = 100000
$ws.Range("N131").Value = -110080

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 75590.5
$ws.Range("J27").Value = 75590.5
$ws.Range("L27").Value = 75590.5
$ws.Range("N27").Value = -75728.5
$ws.Range("H102").Value = 79997.5
$ws.Range("J102").Value = 79997.5
$ws.Range("L102").Value = 79997.5
$ws.Range("N102").Value = -86487.5
$ws.Range("H126").Value = 2418.2
$ws.Range("I126").Value = 2226.2222
$ws.Range("K126").Value = 6678.6666
$ws.Range("M126").Value = -4208.6666
$ws.Range("H127").Value = 58000
$ws.Range("J127").Value = 58000
$ws.Range("L127").Value = 58000
$ws.Range("N127").Value = -67920
$ws.Range("H129").Value = 85000
$ws.Range("J129").Value = 85000
$ws.Range("L129").Value = 85000
$ws.Range("N129").Value = -95000
$ws.Range("H132").Value = 46309880
$ws.Range("I132").Value = 7938079.5
$ws.Range("K132").Value = 23814238.5
$ws.Range("M132").Value = -23811708.5
$ws.Range("H138").Value = 80000
$ws.Range("J138").Value = 80000
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280

Write-Output "Applied all Cactuar_Profits updates"